$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values that look numeric (e.g. "607.33") but are
# stored as literal text in the source data (thousand-separator dots make
# them unsuitable as real numbers, e.g. "69.289.55"). Force text format so
# Excel does not auto-convert/round them to doubles.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.289.55"
$ws.Range("E2").Value = "  +0.65%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.517.79"
$ws.Range("E3").Value = "  +0.53%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.33"
$ws.Range("E5").Value = "  +4.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.89"
$ws.Range("E6").Value = "  -2.91%  "
$ws.Range("E7").Value = "  -1.41%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.512.62"
$ws.Range("E8").Value = "  +0.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.195"
$ws.Range("E10").Value = "  +3.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.65"
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.580"
$ws.Range("E12").Value = "  -3.79%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "47.27"
$ws.Range("E13").Value = "  +0.02%  "
$ws.Range("E14").Value = "  +0.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.095.37"
$ws.Range("E15").Value = "  +0.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.36"
$ws.Range("E16").Value = "  -5.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "614.55"
$ws.Range("E17").Value = "  -9.93%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.523.59"
$ws.Range("E18").Value = "  +0.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.425.50"
$ws.Range("E19").Value = "  +0.87%  "
$ws.Range("E20").Value = "  -1.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.25"
$ws.Range("E21").Value = "  -1.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.22"
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.882"
$ws.Range("E23").Value = "  -2.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.81"
$ws.Range("E24").Value = "  -3.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.41"
$ws.Range("E25").Value = "  -1.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.88"
$ws.Range("E26").Value = "  +1.13%  "
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("E28").Value = "  -0.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.28"
$ws.Range("E29").Value = "  -1.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.12"
$ws.Range("E30").Value = "  +0.42%  "
$ws.Range("E31").Value = "  -2.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.49"
$ws.Range("E32").Value = "  -3.08%  "
$ws.Range("E33").Value = "  -0.70%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.93"
$ws.Range("E34").Value = "  -6.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "567.66"
$ws.Range("E35").Value = "  -0.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.78"
$ws.Range("E36").Value = "  -1.61%  "
$ws.Range("E37").Value = "  -2.91%  "
$ws.Range("E38").Value = "  -3.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "56.95"
$ws.Range("E39").Value = "  +0.22%  "
$ws.Range("E40").Value = "  +0.07%  "
$ws.Range("E41").Value = "  +1.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0444"
$ws.Range("E42").Value = "  +0.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.382.74"
$ws.Range("E43").Value = "  -1.21%  "
$ws.Range("E44").Value = "  -2.83%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "33.06"
$ws.Range("E45").Value = "  -1.31%  "
$ws.Range("E46").Value = "  +0.27%  "
$ws.Range("E47").Value = "  -1.10%  "
$ws.Range("E48").Value = "  +0.14%  "
$ws.Range("E49").Value = "  -3.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.14"
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("E51").Value = "  +8.55%  "
